# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2-73) holds a date serial that was meant to represent the
# "15th of the month" for each quarterly/period marker, but was stored as
# the 1st of the month instead. This shifts every date in column A from the
# 1st of its month to the 15th of the *following* month (matching the
# corrected quarterly indexing), leaving every other cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 73
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2

    if ($serial -eq $null) {
        continue
    }

    $d = [DateTime]::FromOADate($serial)

    $newMonth = $d.Month + 1
    $newYear = $d.Year
    if ($newMonth -gt 12) {
        $newMonth = 1
        $newYear = $newYear + 1
    }

    $newDate = Get-Date -Year $newYear -Month $newMonth -Day 15 -Hour 0 -Minute 0 -Second 0
    $cell.Value = $newDate.ToOADate()
}

Write-Output "Updated column A date serials for rows 2..$lastRow"
